$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 91.20678233333332
$ws.Cells.Item(2, 8).Value = 273.620347
$ws.Cells.Item(2, 9).Value = 0.6532435006323181
$ws.Cells.Item(2, 10).Value = 0.6532435006323182
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 142.158333
$ws.Cells.Item(2, 14).Value = 426.474999
$ws.Cells.Item(2, 15).Value = 0.9500780504241082
$ws.Cells.Item(2, 16).Value = 0.9500780504241081
$ws.Cells.Item(2, 17).Value = 12965.80413480052
$ws.Cells.Item(2, 18).Value = 116692.2372132047
$ws.Cells.Item(2, 19).Value = 0.6206323115329725
$ws.Cells.Item(2, 20).Value = 0.6206323115329725

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 91.20678233333332
$ws.Cells.Item(3, 8).Value = 273.620347
$ws.Cells.Item(3, 9).Value = 0.6532435006323181
$ws.Cells.Item(3, 10).Value = 0.6532435006323182
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.2277323333333333
$ws.Cells.Item(3, 14).Value = 0.6831970000000001
$ws.Cells.Item(3, 15).Value = 0.00152198950779668
$ws.Cells.Item(3, 16).Value = 0.00152198950779668
$ws.Cells.Item(3, 17).Value = 20.77073335659544
$ws.Cells.Item(3, 18).Value = 186.936600209359
$ws.Cells.Item(3, 19).Value = 0.0009942297539987618
$ws.Cells.Item(3, 20).Value = 0.0009942297539987618

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 91.20678233333332
$ws.Cells.Item(4, 8).Value = 273.620347
$ws.Cells.Item(4, 9).Value = 0.6532435006323181
$ws.Cells.Item(4, 10).Value = 0.6532435006323182
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.576418
$ws.Cells.Item(4, 14).Value = 13.729254
$ws.Cells.Item(4, 15).Value = 0.03058529317001626
$ws.Cells.Item(4, 16).Value = 0.03058529317001625
$ws.Cells.Item(4, 17).Value = 417.4003603923487
$ws.Cells.Item(4, 18).Value = 3756.603243531138
$ws.Cells.Item(4, 19).Value = 0.01997964397824715
$ws.Cells.Item(4, 20).Value = 0.01997964397824715

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 91.20678233333332
$ws.Cells.Item(5, 8).Value = 273.620347
$ws.Cells.Item(5, 9).Value = 0.6532435006323181
$ws.Cells.Item(5, 10).Value = 0.6532435006323182
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 2.532751666666667
$ws.Cells.Item(5, 14).Value = 7.598255
$ws.Cells.Item(5, 15).Value = 0.01692698356047181
$ws.Cells.Item(5, 16).Value = 0.01692698356047181
$ws.Cells.Item(5, 17).Value = 231.0041299660539
$ws.Cells.Item(5, 18).Value = 2079.037169694485
$ws.Cells.Item(5, 19).Value = 0.0110574419961883
$ws.Cells.Item(5, 20).Value = 0.0110574419961883

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 91.20678233333332
$ws.Cells.Item(6, 8).Value = 273.620347
$ws.Cells.Item(6, 9).Value = 0.6532435006323181
$ws.Cells.Item(6, 10).Value = 0.6532435006323182
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.1328223333333333
$ws.Cells.Item(6, 14).Value = 0.398467
$ws.Cells.Item(6, 15).Value = 0.0008876833376071904
$ws.Cells.Item(6, 16).Value = 0.00088768333760719
$ws.Cells.Item(6, 17).Value = 12.11429764533878
$ws.Cells.Item(6, 18).Value = 109.028678808049
$ws.Cells.Item(6, 19).Value = 0.0005798733709115009
$ws.Cells.Item(6, 20).Value = 0.0005798733709115008

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 27.97197233333334
$ws.Cells.Item(7, 8).Value = 83.91591700000001
$ws.Cells.Item(7, 9).Value = 0.2003415607825798
$ws.Cells.Item(7, 10).Value = 0.2003415607825798
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 142.158333
$ws.Cells.Item(7, 14).Value = 426.474999
$ws.Cells.Item(7, 15).Value = 0.9500780504241082
$ws.Cells.Item(7, 16).Value = 0.9500780504241081
$ws.Cells.Item(7, 17).Value = 3976.448957628787
$ws.Cells.Item(7, 18).Value = 35788.04061865909
$ws.Cells.Item(7, 19).Value = 0.1903401194872363
$ws.Cells.Item(7, 20).Value = 0.1903401194872363

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 27.97197233333334
$ws.Cells.Item(8, 8).Value = 83.91591700000001
$ws.Cells.Item(8, 9).Value = 0.2003415607825798
$ws.Cells.Item(8, 10).Value = 0.2003415607825798
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.2277323333333333
$ws.Cells.Item(8, 14).Value = 0.6831970000000001
$ws.Cells.Item(8, 15).Value = 0.00152198950779668
$ws.Cells.Item(8, 16).Value = 0.00152198950779668
$ws.Cells.Item(8, 17).Value = 6.370122527405446
$ws.Cells.Item(8, 18).Value = 57.33110274664901
$ws.Cells.Item(8, 19).Value = 0.0003049177534866972
$ws.Cells.Item(8, 20).Value = 0.0003049177534866972

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 27.97197233333334
$ws.Cells.Item(9, 8).Value = 83.91591700000001
$ws.Cells.Item(9, 9).Value = 0.2003415607825798
$ws.Cells.Item(9, 10).Value = 0.2003415607825798
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.576418
$ws.Cells.Item(9, 14).Value = 13.729254
$ws.Cells.Item(9, 15).Value = 0.03058529317001626
$ws.Cells.Item(9, 16).Value = 0.03058529317001625
$ws.Cells.Item(9, 17).Value = 128.0114376817687
$ws.Cells.Item(9, 18).Value = 1152.102939135918
$ws.Cells.Item(9, 19).Value = 0.006127505370673834
$ws.Cells.Item(9, 20).Value = 0.006127505370673833

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 27.97197233333334
$ws.Cells.Item(10, 8).Value = 83.91591700000001
$ws.Cells.Item(10, 9).Value = 0.2003415607825798
$ws.Cells.Item(10, 10).Value = 0.2003415607825798
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.532751666666667
$ws.Cells.Item(10, 14).Value = 7.598255
$ws.Cells.Item(10, 15).Value = 0.01692698356047181
$ws.Cells.Item(10, 16).Value = 0.01692698356047181
$ws.Cells.Item(10, 17).Value = 70.84605954720389
$ws.Cells.Item(10, 18).Value = 637.6145359248351
$ws.Cells.Item(10, 19).Value = 0.003391178305845992
$ws.Cells.Item(10, 20).Value = 0.003391178305845991

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 27.97197233333334
$ws.Cells.Item(11, 8).Value = 83.91591700000001
$ws.Cells.Item(11, 9).Value = 0.2003415607825798
$ws.Cells.Item(11, 10).Value = 0.2003415607825798
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.1328223333333333
$ws.Cells.Item(11, 14).Value = 0.398467
$ws.Cells.Item(11, 15).Value = 0.0008876833376071904
$ws.Cells.Item(11, 16).Value = 0.00088768333760719
$ws.Cells.Item(11, 17).Value = 3.715302633248779
$ws.Cells.Item(11, 18).Value = 33.437723699239
$ws.Cells.Item(11, 19).Value = 0.0001778398653369142
$ws.Cells.Item(11, 20).Value = 0.0001778398653369142

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 20.442661
$ws.Cells.Item(12, 8).Value = 61.327983
$ws.Cells.Item(12, 9).Value = 0.146414938585102
$ws.Cells.Item(12, 10).Value = 0.146414938585102
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 142.158333
$ws.Cells.Item(12, 14).Value = 426.474999
$ws.Cells.Item(12, 15).Value = 0.9500780504241082
$ws.Cells.Item(12, 16).Value = 0.9500780504241081
$ws.Cells.Item(12, 17).Value = 2906.094609844113
$ws.Cells.Item(12, 18).Value = 26154.85148859702
$ws.Cells.Item(12, 19).Value = 0.1391056194038993
$ws.Cells.Item(12, 20).Value = 0.1391056194038992

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 20.442661
$ws.Cells.Item(13, 8).Value = 61.327983
$ws.Cells.Item(13, 9).Value = 0.146414938585102
$ws.Cells.Item(13, 10).Value = 0.146414938585102
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.2277323333333333
$ws.Cells.Item(13, 14).Value = 0.6831970000000001
$ws.Cells.Item(13, 15).Value = 0.00152198950779668
$ws.Cells.Item(13, 16).Value = 0.00152198950779668
$ws.Cells.Item(13, 17).Value = 4.655454889072334
$ws.Cells.Item(13, 18).Value = 41.89909400165101
$ws.Cells.Item(13, 19).Value = 0.0002228420003112205
$ws.Cells.Item(13, 20).Value = 0.0002228420003112205

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 20.442661
$ws.Cells.Item(14, 8).Value = 61.327983
$ws.Cells.Item(14, 9).Value = 0.146414938585102
$ws.Cells.Item(14, 10).Value = 0.146414938585102
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 4.576418
$ws.Cells.Item(14, 14).Value = 13.729254
$ws.Cells.Item(14, 15).Value = 0.03058529317001626
$ws.Cells.Item(14, 16).Value = 0.03058529317001625
$ws.Cells.Item(14, 17).Value = 93.55416176829802
$ws.Cells.Item(14, 18).Value = 841.9874559146821
$ws.Cells.Item(14, 19).Value = 0.00447814382109527
$ws.Cells.Item(14, 20).Value = 0.004478143821095269

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 20.442661
$ws.Cells.Item(15, 8).Value = 61.327983
$ws.Cells.Item(15, 9).Value = 0.146414938585102
$ws.Cells.Item(15, 10).Value = 0.146414938585102
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 2.532751666666667
$ws.Cells.Item(15, 14).Value = 7.598255
$ws.Cells.Item(15, 15).Value = 0.01692698356047181
$ws.Cells.Item(15, 16).Value = 0.01692698356047181
$ws.Cells.Item(15, 17).Value = 51.77618371885167
$ws.Cells.Item(15, 18).Value = 465.985653469665
$ws.Cells.Item(15, 19).Value = 0.002478363258437512
$ws.Cells.Item(15, 20).Value = 0.002478363258437511

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 20.442661
$ws.Cells.Item(16, 8).Value = 61.327983
$ws.Cells.Item(16, 9).Value = 0.146414938585102
$ws.Cells.Item(16, 10).Value = 0.146414938585102
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.1328223333333333
$ws.Cells.Item(16, 14).Value = 0.398467
$ws.Cells.Item(16, 15).Value = 0.0008876833376071904
$ws.Cells.Item(16, 16).Value = 0.00088768333760719
$ws.Cells.Item(16, 17).Value = 2.715241933562334
$ws.Cells.Item(16, 18).Value = 24.437177402061
$ws.Cells.Item(16, 19).Value = 0.0001299701013587752
$ws.Cells.Item(16, 20).Value = 0.0001299701013587751
